# The "Förändrad" (Changed) column C date moves forward one day,
# from 2023-10-03 (serial 45202) to 2023-10-04 (serial 45203),
# for every data row (rows 2 through 45) on the active sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 45; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    $val = $cell.Value()
    if ($val.Year() -eq 2023 -and $val.Month() -eq 10 -and $val.Day() -eq 3) {
        $cell.Value = 45203
    }
}
